{"js": "// Add a new paragraph containing \"balance\" right after the existing\n// \"tes cp\" paragraph (and before the section break at the end of the body).\nconst body = context.document.body;\n\n// Locate the \"tes cp\" paragraph so the new paragraph is anchored to it\n// (rather than simply assuming it is the last paragraph in the body).\nconst results = body.search(\"tes cp\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  // Insert a new paragraph right after the found text; it inherits the\n  // same run formatting (lang=\"en-US\") already present on that paragraph.\n  target.insertParagraph(\"balance\", Word.InsertLocation.after);\n} else {\n  // Fallback: if the text could not be found, just append to the body.\n  body.insertParagraph(\"balance\", Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# Add a new paragraph containing \"balance\" right after the existing\n# \"tes cp\" paragraph (and before the section break at the end of the body).\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"tes cp\")\n\nif ($found) {\n    # Collapse the found range to its end (right after \"tes cp\") and split\n    # off a new paragraph there.\n    $rng.Collapse(0)\n    $rng.InsertParagraphAfter()\n\n    # The newly inserted paragraph is now the last paragraph in the\n    # document; give it the \"balance\" text (it already inherited the\n    # lang=\"en-US\" run formatting from the \"tes cp\" paragraph).\n    $newPara = $d.Paragraphs.Last\n    $newPara.Range.Text = \"balance\"\n}\n"}
